$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 244
$ws.Cells.Item(2, 7).Value = "bedrooms"
$ws.Cells.Item(2, 8).Value = "living_rooms"
$ws.Cells.Item(2, 9).Value = "distractor"
$ws.Cells.Item(2, 11).Value = "f"
$ws.Cells.Item(2, 12).Value = "stimuli/img_hc49v.png"
$ws.Cells.Item(2, 13).Value = 70.95121951219512
$ws.Cells.Item(2, 14).Value = 53.31707317073171
$ws.Cells.Item(2, 15).Value = 62.13414634146342
$ws.Cells.Item(2, 16).Value = 41
$ws.Cells.Item(2, 17).Value = 6
$ws.Cells.Item(2, 18).Value = 6
$ws.Cells.Item(2, 19).Value = 6
$ws.Cells.Item(2, 20).Value = 6
$ws.Cells.Item(2, 21).Value = 6
$ws.Cells.Item(2, 22).Value = 6

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 245
$ws.Cells.Item(3, 7).Value = "bedrooms"
$ws.Cells.Item(3, 8).Value = "bedrooms"
$ws.Cells.Item(3, 9).Value = "target"
$ws.Cells.Item(3, 11).Value = "j"
$ws.Cells.Item(3, 12).Value = "stimuli/img_ose78.png"
$ws.Cells.Item(3, 13).Value = 80.19444444444444
$ws.Cells.Item(3, 14).Value = 60.25
$ws.Cells.Item(3, 15).Value = 70.22222222222223
$ws.Cells.Item(3, 16).Value = 36
$ws.Cells.Item(3, 17).Value = 8
$ws.Cells.Item(3, 18).Value = 7
$ws.Cells.Item(3, 19).Value = 7
$ws.Cells.Item(3, 20).Value = 7
$ws.Cells.Item(3, 21).Value = 7
$ws.Cells.Item(3, 22).Value = 7

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 246
$ws.Cells.Item(4, 7).Value = "bedrooms"
$ws.Cells.Item(4, 8).Value = "kitchens"
$ws.Cells.Item(4, 9).Value = "distractor"
$ws.Cells.Item(4, 11).Value = "f"
$ws.Cells.Item(4, 12).Value = "stimuli/img_oau79.png"
$ws.Cells.Item(4, 13).Value = 70.86486486486487
$ws.Cells.Item(4, 14).Value = 49
$ws.Cells.Item(4, 15).Value = 59.93243243243244
$ws.Cells.Item(4, 16).Value = 37
$ws.Cells.Item(4, 17).Value = 5
$ws.Cells.Item(4, 18).Value = 5
$ws.Cells.Item(4, 19).Value = 5
$ws.Cells.Item(4, 20).Value = 5
$ws.Cells.Item(4, 21).Value = 5
$ws.Cells.Item(4, 22).Value = 5

$ws.Cells.Item(5, 5).Value = 4
$ws.Cells.Item(5, 6).Value = 247
$ws.Cells.Item(5, 7).Value = "bedrooms"
$ws.Cells.Item(5, 8).Value = "living_rooms"
$ws.Cells.Item(5, 9).Value = "distractor"
$ws.Cells.Item(5, 11).Value = "f"
$ws.Cells.Item(5, 12).Value = "stimuli/img_53nbn.png"
$ws.Cells.Item(5, 13).Value = 73.28888888888889
$ws.Cells.Item(5, 14).Value = 51.15555555555556
$ws.Cells.Item(5, 15).Value = 62.22222222222223
$ws.Cells.Item(5, 16).Value = 45
$ws.Cells.Item(5, 17).Value = 6
$ws.Cells.Item(5, 18).Value = 6
$ws.Cells.Item(5, 19).Value = 6
$ws.Cells.Item(5, 20).Value = 6
$ws.Cells.Item(5, 21).Value = 6
$ws.Cells.Item(5, 22).Value = 6

$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = 248
$ws.Cells.Item(6, 7).Value = "bedrooms"
$ws.Cells.Item(6, 8).Value = "living_rooms"
$ws.Cells.Item(6, 9).Value = "distractor"
$ws.Cells.Item(6, 11).Value = "f"
$ws.Cells.Item(6, 12).Value = "stimuli/img_wz6x5.png"
$ws.Cells.Item(6, 13).Value = 68.3695652173913
$ws.Cells.Item(6, 14).Value = 48.47826086956522
$ws.Cells.Item(6, 15).Value = 58.42391304347826
$ws.Cells.Item(6, 16).Value = 46
$ws.Cells.Item(6, 17).Value = 5
$ws.Cells.Item(6, 18).Value = 5
$ws.Cells.Item(6, 19).Value = 5
$ws.Cells.Item(6, 20).Value = 5
$ws.Cells.Item(6, 21).Value = 5
$ws.Cells.Item(6, 22).Value = 5

$ws.Cells.Item(7, 5).Value = 6
$ws.Cells.Item(7, 6).Value = 249
$ws.Cells.Item(7, 7).Value = "bedrooms"
$ws.Cells.Item(7, 8).Value = "bedrooms"
$ws.Cells.Item(7, 9).Value = "target"
$ws.Cells.Item(7, 11).Value = "j"
$ws.Cells.Item(7, 12).Value = "stimuli/img_5p2ql.png"
$ws.Cells.Item(7, 13).Value = 89.19565217391305
$ws.Cells.Item(7, 14).Value = 72.52173913043478
$ws.Cells.Item(7, 15).Value = 80.8586956521739
$ws.Cells.Item(7, 16).Value = 46
$ws.Cells.Item(7, 17).Value = 10
$ws.Cells.Item(7, 18).Value = 10
$ws.Cells.Item(7, 19).Value = 10
$ws.Cells.Item(7, 20).Value = 10
$ws.Cells.Item(7, 21).Value = 10
$ws.Cells.Item(7, 22).Value = 9

$ws.Cells.Item(8, 5).Value = 7
$ws.Cells.Item(8, 6).Value = 250
$ws.Cells.Item(8, 7).Value = "bedrooms"
$ws.Cells.Item(8, 8).Value = "kitchens"
$ws.Cells.Item(8, 9).Value = "distractor"
$ws.Cells.Item(8, 11).Value = "f"
$ws.Cells.Item(8, 12).Value = "stimuli/img_68wfw.png"
$ws.Cells.Item(8, 13).Value = 68.87878787878788
$ws.Cells.Item(8, 14).Value = 48.96969696969697
$ws.Cells.Item(8, 15).Value = 58.92424242424242
$ws.Cells.Item(8, 16).Value = 33
$ws.Cells.Item(8, 17).Value = 5
$ws.Cells.Item(8, 18).Value = 5
$ws.Cells.Item(8, 19).Value = 5
$ws.Cells.Item(8, 20).Value = 5
$ws.Cells.Item(8, 21).Value = 5
$ws.Cells.Item(8, 22).Value = 5

$ws.Cells.Item(9, 5).Value = 8
$ws.Cells.Item(9, 6).Value = 251
$ws.Cells.Item(9, 7).Value = "bedrooms"
$ws.Cells.Item(9, 8).Value = "bedrooms"
$ws.Cells.Item(9, 9).Value = "target"
$ws.Cells.Item(9, 11).Value = "j"
$ws.Cells.Item(9, 12).Value = "stimuli/img_gbypq.png"
$ws.Cells.Item(9, 13).Value = 76.275
$ws.Cells.Item(9, 14).Value = 51.925
$ws.Cells.Item(9, 15).Value = 64.1
$ws.Cells.Item(9, 16).Value = 40
$ws.Cells.Item(9, 17).Value = 6
$ws.Cells.Item(9, 18).Value = 6
$ws.Cells.Item(9, 19).Value = 6
$ws.Cells.Item(9, 20).Value = 6
$ws.Cells.Item(9, 21).Value = 6
$ws.Cells.Item(9, 22).Value = 6

$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 252
$ws.Cells.Item(10, 7).Value = "bedrooms"
$ws.Cells.Item(10, 8).Value = "kitchens"
$ws.Cells.Item(10, 9).Value = "distractor"
$ws.Cells.Item(10, 11).Value = "f"
$ws.Cells.Item(10, 12).Value = "stimuli/img_mucwi.png"
$ws.Cells.Item(10, 13).Value = 71.14814814814815
$ws.Cells.Item(10, 14).Value = 48.55555555555556
$ws.Cells.Item(10, 15).Value = 59.85185185185185
$ws.Cells.Item(10, 16).Value = 27
$ws.Cells.Item(10, 17).Value = 5
$ws.Cells.Item(10, 18).Value = 5
$ws.Cells.Item(10, 19).Value = 5
$ws.Cells.Item(10, 20).Value = 5
$ws.Cells.Item(10, 21).Value = 5
$ws.Cells.Item(10, 22).Value = 5

$ws.Cells.Item(11, 5).Value = 10
$ws.Cells.Item(11, 6).Value = 253
$ws.Cells.Item(11, 7).Value = "bedrooms"
$ws.Cells.Item(11, 8).Value = "bedrooms"
$ws.Cells.Item(11, 9).Value = "target"
$ws.Cells.Item(11, 11).Value = "j"
$ws.Cells.Item(11, 12).Value = "stimuli/img_aweye.png"
$ws.Cells.Item(11, 13).Value = 53.42105263157895
$ws.Cells.Item(11, 14).Value = 31.84210526315789
$ws.Cells.Item(11, 15).Value = 42.63157894736842
$ws.Cells.Item(11, 16).Value = 38
$ws.Cells.Item(11, 17).Value = 2
$ws.Cells.Item(11, 18).Value = 2
$ws.Cells.Item(11, 19).Value = 2
$ws.Cells.Item(11, 20).Value = 3
$ws.Cells.Item(11, 21).Value = 3
$ws.Cells.Item(11, 22).Value = 2

$ws.Cells.Item(12, 5).Value = 11
$ws.Cells.Item(12, 6).Value = 254
$ws.Cells.Item(12, 7).Value = "bedrooms"
$ws.Cells.Item(12, 8).Value = "living_rooms"
$ws.Cells.Item(12, 9).Value = "distractor"
$ws.Cells.Item(12, 11).Value = "f"
$ws.Cells.Item(12, 12).Value = "stimuli/img_lgxzn.png"
$ws.Cells.Item(12, 13).Value = 73.11363636363636
$ws.Cells.Item(12, 14).Value = 49.97727272727273
$ws.Cells.Item(12, 15).Value = 61.54545454545455
$ws.Cells.Item(12, 16).Value = 44
$ws.Cells.Item(12, 17).Value = 6
$ws.Cells.Item(12, 18).Value = 6
$ws.Cells.Item(12, 19).Value = 6
$ws.Cells.Item(12, 20).Value = 6
$ws.Cells.Item(12, 21).Value = 6
$ws.Cells.Item(12, 22).Value = 5

$ws.Cells.Item(13, 5).Value = 12
$ws.Cells.Item(13, 6).Value = 255
$ws.Cells.Item(13, 7).Value = "bedrooms"
$ws.Cells.Item(13, 8).Value = "bedrooms"
$ws.Cells.Item(13, 9).Value = "target"
$ws.Cells.Item(13, 11).Value = "j"
$ws.Cells.Item(13, 12).Value = "stimuli/img_bj2gr.png"
$ws.Cells.Item(13, 13).Value = 65.25
$ws.Cells.Item(13, 14).Value = 44.8
$ws.Cells.Item(13, 15).Value = 55.025
$ws.Cells.Item(13, 16).Value = 40
$ws.Cells.Item(13, 17).Value = 4
$ws.Cells.Item(13, 18).Value = 4
$ws.Cells.Item(13, 19).Value = 4
$ws.Cells.Item(13, 20).Value = 4
$ws.Cells.Item(13, 21).Value = 4
$ws.Cells.Item(13, 22).Value = 4

$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = 256
$ws.Cells.Item(14, 7).Value = "bedrooms"
$ws.Cells.Item(14, 8).Value = "bedrooms"
$ws.Cells.Item(14, 9).Value = "target"
$ws.Cells.Item(14, 11).Value = "j"
$ws.Cells.Item(14, 12).Value = "stimuli/img_okvvw.png"
$ws.Cells.Item(14, 13).Value = 50.58333333333334
$ws.Cells.Item(14, 14).Value = 32.11111111111111
$ws.Cells.Item(14, 15).Value = 41.34722222222223
$ws.Cells.Item(14, 16).Value = 36
$ws.Cells.Item(14, 17).Value = 2
$ws.Cells.Item(14, 18).Value = 2
$ws.Cells.Item(14, 19).Value = 2
$ws.Cells.Item(14, 20).Value = 2
$ws.Cells.Item(14, 21).Value = 2
$ws.Cells.Item(14, 22).Value = 3

$ws.Cells.Item(15, 5).Value = 14
$ws.Cells.Item(15, 6).Value = 257
$ws.Cells.Item(15, 7).Value = "bedrooms"
$ws.Cells.Item(15, 8).Value = "kitchens"
$ws.Cells.Item(15, 9).Value = "distractor"
$ws.Cells.Item(15, 11).Value = "f"
$ws.Cells.Item(15, 12).Value = "stimuli/img_d0k76.png"
$ws.Cells.Item(15, 13).Value = 67.0909090909091
$ws.Cells.Item(15, 14).Value = 46.3030303030303
$ws.Cells.Item(15, 15).Value = 56.6969696969697
$ws.Cells.Item(15, 16).Value = 33
$ws.Cells.Item(15, 17).Value = 4
$ws.Cells.Item(15, 18).Value = 4
$ws.Cells.Item(15, 19).Value = 4
$ws.Cells.Item(15, 20).Value = 4
$ws.Cells.Item(15, 21).Value = 4
$ws.Cells.Item(15, 22).Value = 4

$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 258
$ws.Cells.Item(16, 7).Value = "bedrooms"
$ws.Cells.Item(16, 8).Value = "kitchens"
$ws.Cells.Item(16, 9).Value = "distractor"
$ws.Cells.Item(16, 11).Value = "f"
$ws.Cells.Item(16, 12).Value = "stimuli/img_pt3d7.png"
$ws.Cells.Item(16, 13).Value = 65.08571428571429
$ws.Cells.Item(16, 14).Value = 44.65714285714286
$ws.Cells.Item(16, 15).Value = 54.87142857142857
$ws.Cells.Item(16, 16).Value = 35
$ws.Cells.Item(16, 17).Value = 4
$ws.Cells.Item(16, 18).Value = 4
$ws.Cells.Item(16, 19).Value = 4
$ws.Cells.Item(16, 20).Value = 4
$ws.Cells.Item(16, 21).Value = 4
$ws.Cells.Item(16, 22).Value = 4

$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 6).Value = 259
$ws.Cells.Item(17, 7).Value = "bedrooms"
$ws.Cells.Item(17, 8).Value = "bedrooms"
$ws.Cells.Item(17, 9).Value = "target"
$ws.Cells.Item(17, 11).Value = "j"
$ws.Cells.Item(17, 12).Value = "stimuli/img_v8dra.png"
$ws.Cells.Item(17, 13).Value = 61.77272727272727
$ws.Cells.Item(17, 14).Value = 38.79545454545455
$ws.Cells.Item(17, 15).Value = 50.28409090909091
$ws.Cells.Item(17, 16).Value = 44
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = 3
$ws.Cells.Item(17, 19).Value = 3
$ws.Cells.Item(17, 20).Value = 3
$ws.Cells.Item(17, 21).Value = 4
$ws.Cells.Item(17, 22).Value = 3

$ws.Cells.Item(18, 5).Value = 17
$ws.Cells.Item(18, 6).Value = 260
$ws.Cells.Item(18, 7).Value = "bedrooms"
$ws.Cells.Item(18, 8).Value = "kitchens"
$ws.Cells.Item(18, 9).Value = "distractor"
$ws.Cells.Item(18, 11).Value = "f"
$ws.Cells.Item(18, 12).Value = "stimuli/img_g7870.png"
$ws.Cells.Item(18, 13).Value = 68.70967741935483
$ws.Cells.Item(18, 14).Value = 44.2258064516129
$ws.Cells.Item(18, 15).Value = 56.46774193548387
$ws.Cells.Item(18, 16).Value = 31
$ws.Cells.Item(18, 17).Value = 4
$ws.Cells.Item(18, 18).Value = 4
$ws.Cells.Item(18, 19).Value = 4
$ws.Cells.Item(18, 20).Value = 4
$ws.Cells.Item(18, 21).Value = 4
$ws.Cells.Item(18, 22).Value = 4

$ws.Cells.Item(19, 5).Value = 18
$ws.Cells.Item(19, 6).Value = 261
$ws.Cells.Item(19, 7).Value = "bedrooms"
$ws.Cells.Item(19, 8).Value = "bedrooms"
$ws.Cells.Item(19, 9).Value = "target"
$ws.Cells.Item(19, 11).Value = "j"
$ws.Cells.Item(19, 12).Value = "stimuli/img_t4hvr.png"
$ws.Cells.Item(19, 13).Value = 61.69230769230769
$ws.Cells.Item(19, 14).Value = 39.76923076923077
$ws.Cells.Item(19, 15).Value = 50.73076923076923
$ws.Cells.Item(19, 16).Value = 39
$ws.Cells.Item(19, 17).Value = 3
$ws.Cells.Item(19, 18).Value = 3
$ws.Cells.Item(19, 19).Value = 3
$ws.Cells.Item(19, 20).Value = 4
$ws.Cells.Item(19, 21).Value = 3
$ws.Cells.Item(19, 22).Value = 4

$ws.Cells.Item(20, 5).Value = 19
$ws.Cells.Item(20, 6).Value = 262
$ws.Cells.Item(20, 7).Value = "bedrooms"
$ws.Cells.Item(20, 8).Value = "bedrooms"
$ws.Cells.Item(20, 9).Value = "target"
$ws.Cells.Item(20, 11).Value = "j"
$ws.Cells.Item(20, 12).Value = "stimuli/img_kzg3h.png"
$ws.Cells.Item(20, 13).Value = 77.02777777777777
$ws.Cells.Item(20, 14).Value = 56.22222222222222
$ws.Cells.Item(20, 15).Value = 66.625
$ws.Cells.Item(20, 16).Value = 36
$ws.Cells.Item(20, 17).Value = 7
$ws.Cells.Item(20, 18).Value = 7
$ws.Cells.Item(20, 19).Value = 7
$ws.Cells.Item(20, 20).Value = 7
$ws.Cells.Item(20, 21).Value = 7
$ws.Cells.Item(20, 22).Value = 7

$ws.Cells.Item(21, 5).Value = 20
$ws.Cells.Item(21, 6).Value = 263
$ws.Cells.Item(21, 7).Value = "bedrooms"
$ws.Cells.Item(21, 8).Value = "bedrooms"
$ws.Cells.Item(21, 9).Value = "target"
$ws.Cells.Item(21, 11).Value = "j"
$ws.Cells.Item(21, 12).Value = "stimuli/img_2pk6v.png"
$ws.Cells.Item(21, 13).Value = 85.08108108108108
$ws.Cells.Item(21, 14).Value = 66.16216216216216
$ws.Cells.Item(21, 15).Value = 75.62162162162161
$ws.Cells.Item(21, 16).Value = 37
$ws.Cells.Item(21, 17).Value = 9
$ws.Cells.Item(21, 18).Value = 9
$ws.Cells.Item(21, 19).Value = 9
$ws.Cells.Item(21, 20).Value = 9
$ws.Cells.Item(21, 21).Value = 9
$ws.Cells.Item(21, 22).Value = 8

$ws.Cells.Item(22, 5).Value = 21
$ws.Cells.Item(22, 6).Value = 264
$ws.Cells.Item(22, 7).Value = "bedrooms"
$ws.Cells.Item(22, 8).Value = "kitchens"
$ws.Cells.Item(22, 9).Value = "distractor"
$ws.Cells.Item(22, 11).Value = "f"
$ws.Cells.Item(22, 12).Value = "stimuli/img_b971s.png"
$ws.Cells.Item(22, 13).Value = 70.5
$ws.Cells.Item(22, 14).Value = 47.61111111111111
$ws.Cells.Item(22, 15).Value = 59.05555555555556
$ws.Cells.Item(22, 16).Value = 36
$ws.Cells.Item(22, 17).Value = 5
$ws.Cells.Item(22, 18).Value = 5
$ws.Cells.Item(22, 19).Value = 5
$ws.Cells.Item(22, 20).Value = 5
$ws.Cells.Item(22, 21).Value = 5
$ws.Cells.Item(22, 22).Value = 5

$ws.Cells.Item(23, 5).Value = 22
$ws.Cells.Item(23, 6).Value = 265
$ws.Cells.Item(23, 7).Value = "bedrooms"
$ws.Cells.Item(23, 8).Value = "living_rooms"
$ws.Cells.Item(23, 9).Value = "distractor"
$ws.Cells.Item(23, 11).Value = "f"
$ws.Cells.Item(23, 12).Value = "stimuli/img_73pyk.png"
$ws.Cells.Item(23, 13).Value = 69.27659574468085
$ws.Cells.Item(23, 14).Value = 47.27659574468085
$ws.Cells.Item(23, 15).Value = 58.27659574468085
$ws.Cells.Item(23, 16).Value = 47
$ws.Cells.Item(23, 17).Value = 5
$ws.Cells.Item(23, 18).Value = 5
$ws.Cells.Item(23, 19).Value = 5
$ws.Cells.Item(23, 20).Value = 5
$ws.Cells.Item(23, 21).Value = 5
$ws.Cells.Item(23, 22).Value = 5

$ws.Cells.Item(24, 5).Value = 23
$ws.Cells.Item(24, 6).Value = 266
$ws.Cells.Item(24, 7).Value = "bedrooms"
$ws.Cells.Item(24, 8).Value = "living_rooms"
$ws.Cells.Item(24, 9).Value = "distractor"
$ws.Cells.Item(24, 11).Value = "f"
$ws.Cells.Item(24, 12).Value = "stimuli/img_swq34.png"
$ws.Cells.Item(24, 13).Value = 64.11363636363636
$ws.Cells.Item(24, 14).Value = 43.04545454545455
$ws.Cells.Item(24, 15).Value = 53.57954545454545
$ws.Cells.Item(24, 16).Value = 44
$ws.Cells.Item(24, 17).Value = 5
$ws.Cells.Item(24, 18).Value = 5
$ws.Cells.Item(24, 19).Value = 5
$ws.Cells.Item(24, 20).Value = 5
$ws.Cells.Item(24, 21).Value = 5
$ws.Cells.Item(24, 22).Value = 5

$ws.Cells.Item(25, 5).Value = 24
$ws.Cells.Item(25, 6).Value = 267
$ws.Cells.Item(25, 7).Value = "bedrooms"
$ws.Cells.Item(25, 8).Value = "kitchens"
$ws.Cells.Item(25, 9).Value = "distractor"
$ws.Cells.Item(25, 11).Value = "f"
$ws.Cells.Item(25, 12).Value = "stimuli/img_q1ynd.png"
$ws.Cells.Item(25, 13).Value = 70.05714285714286
$ws.Cells.Item(25, 14).Value = 47.31428571428572
$ws.Cells.Item(25, 15).Value = 58.68571428571429
$ws.Cells.Item(25, 16).Value = 35
$ws.Cells.Item(25, 17).Value = 5
$ws.Cells.Item(25, 18).Value = 5
$ws.Cells.Item(25, 19).Value = 5
$ws.Cells.Item(25, 20).Value = 5
$ws.Cells.Item(25, 21).Value = 5
$ws.Cells.Item(25, 22).Value = 5

$ws.Cells.Item(26, 5).Value = 25
$ws.Cells.Item(26, 6).Value = 268
$ws.Cells.Item(26, 7).Value = "bedrooms"
$ws.Cells.Item(26, 8).Value = "kitchens"
$ws.Cells.Item(26, 9).Value = "distractor"
$ws.Cells.Item(26, 11).Value = "f"
$ws.Cells.Item(26, 12).Value = "stimuli/img_4ufga.png"
$ws.Cells.Item(26, 13).Value = 67.79411764705883
$ws.Cells.Item(26, 14).Value = 41.5
$ws.Cells.Item(26, 15).Value = 54.64705882352941
$ws.Cells.Item(26, 16).Value = 34
$ws.Cells.Item(26, 17).Value = 4
$ws.Cells.Item(26, 18).Value = 4
$ws.Cells.Item(26, 19).Value = 4
$ws.Cells.Item(26, 20).Value = 4
$ws.Cells.Item(26, 21).Value = 4
$ws.Cells.Item(26, 22).Value = 4

$ws.Cells.Item(27, 5).Value = 26
$ws.Cells.Item(27, 6).Value = 269
$ws.Cells.Item(27, 7).Value = "bedrooms"
$ws.Cells.Item(27, 8).Value = "bedrooms"
$ws.Cells.Item(27, 9).Value = "target"
$ws.Cells.Item(27, 11).Value = "j"
$ws.Cells.Item(27, 12).Value = "stimuli/img_ic3os.png"
$ws.Cells.Item(27, 13).Value = 84.79069767441861
$ws.Cells.Item(27, 14).Value = 66.16279069767442
$ws.Cells.Item(27, 15).Value = 75.47674418604652
$ws.Cells.Item(27, 16).Value = 43
$ws.Cells.Item(27, 17).Value = 9
$ws.Cells.Item(27, 18).Value = 9
$ws.Cells.Item(27, 19).Value = 9
$ws.Cells.Item(27, 20).Value = 8
$ws.Cells.Item(27, 21).Value = 9
$ws.Cells.Item(27, 22).Value = 9

$ws.Cells.Item(28, 5).Value = 27
$ws.Cells.Item(28, 6).Value = 270
$ws.Cells.Item(28, 7).Value = "bedrooms"
$ws.Cells.Item(28, 8).Value = "bedrooms"
$ws.Cells.Item(28, 9).Value = "target"
$ws.Cells.Item(28, 11).Value = "j"
$ws.Cells.Item(28, 12).Value = "stimuli/img_z3yzz.png"
$ws.Cells.Item(28, 13).Value = 71.71052631578948
$ws.Cells.Item(28, 14).Value = 49.81578947368421
$ws.Cells.Item(28, 15).Value = 60.76315789473685
$ws.Cells.Item(28, 16).Value = 38
$ws.Cells.Item(28, 17).Value = 5
$ws.Cells.Item(28, 18).Value = 5
$ws.Cells.Item(28, 19).Value = 5
$ws.Cells.Item(28, 20).Value = 5
$ws.Cells.Item(28, 21).Value = 5
$ws.Cells.Item(28, 22).Value = 5

$ws.Cells.Item(29, 5).Value = 28
$ws.Cells.Item(29, 6).Value = 271
$ws.Cells.Item(29, 7).Value = "bedrooms"
$ws.Cells.Item(29, 8).Value = "living_rooms"
$ws.Cells.Item(29, 9).Value = "distractor"
$ws.Cells.Item(29, 11).Value = "f"
$ws.Cells.Item(29, 12).Value = "stimuli/img_koooi.png"
$ws.Cells.Item(29, 13).Value = 63.95454545454545
$ws.Cells.Item(29, 14).Value = 44.56818181818182
$ws.Cells.Item(29, 15).Value = 54.26136363636364
$ws.Cells.Item(29, 16).Value = 44
$ws.Cells.Item(29, 17).Value = 5
$ws.Cells.Item(29, 18).Value = 5
$ws.Cells.Item(29, 19).Value = 5
$ws.Cells.Item(29, 20).Value = 5
$ws.Cells.Item(29, 21).Value = 5
$ws.Cells.Item(29, 22).Value = 5

$ws.Cells.Item(30, 5).Value = 29
$ws.Cells.Item(30, 6).Value = 272
$ws.Cells.Item(30, 7).Value = "bedrooms"
$ws.Cells.Item(30, 8).Value = "kitchens"
$ws.Cells.Item(30, 9).Value = "distractor"
$ws.Cells.Item(30, 11).Value = "f"
$ws.Cells.Item(30, 12).Value = "stimuli/img_anjr0.png"
$ws.Cells.Item(30, 13).Value = 67.88888888888889
$ws.Cells.Item(30, 14).Value = 45.80555555555556
$ws.Cells.Item(30, 15).Value = 56.84722222222222
$ws.Cells.Item(30, 16).Value = 36
$ws.Cells.Item(30, 17).Value = 4
$ws.Cells.Item(30, 18).Value = 4
$ws.Cells.Item(30, 19).Value = 4
$ws.Cells.Item(30, 20).Value = 4
$ws.Cells.Item(30, 21).Value = 4
$ws.Cells.Item(30, 22).Value = 4

$ws.Cells.Item(31, 5).Value = 30
$ws.Cells.Item(31, 6).Value = 273
$ws.Cells.Item(31, 7).Value = "bedrooms"
$ws.Cells.Item(31, 8).Value = "bedrooms"
$ws.Cells.Item(31, 9).Value = "target"
$ws.Cells.Item(31, 11).Value = "j"
$ws.Cells.Item(31, 12).Value = "stimuli/img_9pfbj.png"
$ws.Cells.Item(31, 13).Value = 91.27272727272727
$ws.Cells.Item(31, 14).Value = 80.0909090909091
$ws.Cells.Item(31, 15).Value = 85.68181818181819
$ws.Cells.Item(31, 16).Value = 33
$ws.Cells.Item(31, 17).Value = 10
$ws.Cells.Item(31, 18).Value = 10
$ws.Cells.Item(31, 19).Value = 10
$ws.Cells.Item(31, 20).Value = 10
$ws.Cells.Item(31, 21).Value = 10
$ws.Cells.Item(31, 22).Value = 10

$ws.Cells.Item(32, 5).Value = 31
$ws.Cells.Item(32, 6).Value = 274
$ws.Cells.Item(32, 7).Value = "bedrooms"
$ws.Cells.Item(32, 8).Value = "bedrooms"
$ws.Cells.Item(32, 9).Value = "target"
$ws.Cells.Item(32, 11).Value = "j"
$ws.Cells.Item(32, 12).Value = "stimuli/img_2pnl2.png"
$ws.Cells.Item(32, 13).Value = 6.621621621621622
$ws.Cells.Item(32, 14).Value = 7.135135135135135
$ws.Cells.Item(32, 15).Value = 6.878378378378379
$ws.Cells.Item(32, 16).Value = 37
$ws.Cells.Item(32, 17).Value = 1
$ws.Cells.Item(32, 18).Value = 1
$ws.Cells.Item(32, 19).Value = 1
$ws.Cells.Item(32, 20).Value = 1
$ws.Cells.Item(32, 21).Value = 1
$ws.Cells.Item(32, 22).Value = 1

$ws.Cells.Item(33, 5).Value = 32
$ws.Cells.Item(33, 6).Value = 275
$ws.Cells.Item(33, 7).Value = "bedrooms"
$ws.Cells.Item(33, 8).Value = "bedrooms"
$ws.Cells.Item(33, 9).Value = "target"
$ws.Cells.Item(33, 11).Value = "j"
$ws.Cells.Item(33, 12).Value = "stimuli/img_anzgh.png"
$ws.Cells.Item(33, 13).Value = 75.10526315789474
$ws.Cells.Item(33, 14).Value = 55.76315789473684
$ws.Cells.Item(33, 15).Value = 65.4342105263158
$ws.Cells.Item(33, 16).Value = 38
$ws.Cells.Item(33, 17).Value = 6
$ws.Cells.Item(33, 18).Value = 6
$ws.Cells.Item(33, 19).Value = 6
$ws.Cells.Item(33, 20).Value = 6
$ws.Cells.Item(33, 21).Value = 6
$ws.Cells.Item(33, 22).Value = 6

$ws.Cells.Item(34, 5).Value = 33
$ws.Cells.Item(34, 6).Value = 276
$ws.Cells.Item(34, 7).Value = "bedrooms"
$ws.Cells.Item(34, 8).Value = "living_rooms"
$ws.Cells.Item(34, 9).Value = "distractor"
$ws.Cells.Item(34, 11).Value = "f"
$ws.Cells.Item(34, 12).Value = "stimuli/img_xy930.png"
$ws.Cells.Item(34, 13).Value = 70.5952380952381
$ws.Cells.Item(34, 14).Value = 49.47619047619047
$ws.Cells.Item(34, 15).Value = 60.03571428571429
$ws.Cells.Item(34, 16).Value = 42
$ws.Cells.Item(34, 17).Value = 6
$ws.Cells.Item(34, 18).Value = 6
$ws.Cells.Item(34, 19).Value = 6
$ws.Cells.Item(34, 20).Value = 5
$ws.Cells.Item(34, 21).Value = 5
$ws.Cells.Item(34, 22).Value = 5

$ws.Cells.Item(35, 5).Value = 34
$ws.Cells.Item(35, 6).Value = 277
$ws.Cells.Item(35, 7).Value = "bedrooms"
$ws.Cells.Item(35, 8).Value = "bedrooms"
$ws.Cells.Item(35, 9).Value = "target"
$ws.Cells.Item(35, 11).Value = "j"
$ws.Cells.Item(35, 12).Value = "stimuli/img_jivhq.png"
$ws.Cells.Item(35, 13).Value = 37
$ws.Cells.Item(35, 14).Value = 22.26530612244898
$ws.Cells.Item(35, 15).Value = 29.63265306122449
$ws.Cells.Item(35, 16).Value = 49
$ws.Cells.Item(35, 17).Value = 2
$ws.Cells.Item(35, 18).Value = 2
$ws.Cells.Item(35, 19).Value = 2
$ws.Cells.Item(35, 20).Value = 2
$ws.Cells.Item(35, 21).Value = 2
$ws.Cells.Item(35, 22).Value = 2

$ws.Cells.Item(36, 5).Value = 35
$ws.Cells.Item(36, 6).Value = 278
$ws.Cells.Item(36, 7).Value = "bedrooms"
$ws.Cells.Item(36, 8).Value = "bedrooms"
$ws.Cells.Item(36, 9).Value = "target"
$ws.Cells.Item(36, 11).Value = "j"
$ws.Cells.Item(36, 12).Value = "stimuli/img_fqgem.png"
$ws.Cells.Item(36, 13).Value = 80.75
$ws.Cells.Item(36, 14).Value = 61.475
$ws.Cells.Item(36, 15).Value = 71.1125
$ws.Cells.Item(36, 16).Value = 40
$ws.Cells.Item(36, 17).Value = 8
$ws.Cells.Item(36, 18).Value = 8
$ws.Cells.Item(36, 19).Value = 8
$ws.Cells.Item(36, 20).Value = 8
$ws.Cells.Item(36, 21).Value = 8
$ws.Cells.Item(36, 22).Value = 8

$ws.Cells.Item(37, 5).Value = 36
$ws.Cells.Item(37, 6).Value = 279
$ws.Cells.Item(37, 7).Value = "bedrooms"
$ws.Cells.Item(37, 8).Value = "bedrooms"
$ws.Cells.Item(37, 9).Value = "target"
$ws.Cells.Item(37, 11).Value = "j"
$ws.Cells.Item(37, 12).Value = "stimuli/img_yteqw.png"
$ws.Cells.Item(37, 13).Value = 66.83783783783784
$ws.Cells.Item(37, 14).Value = 43.78378378378378
$ws.Cells.Item(37, 15).Value = 55.31081081081081
$ws.Cells.Item(37, 16).Value = 37
$ws.Cells.Item(37, 17).Value = 4
$ws.Cells.Item(37, 18).Value = 4
$ws.Cells.Item(37, 19).Value = 4
$ws.Cells.Item(37, 20).Value = 5
$ws.Cells.Item(37, 21).Value = 4
$ws.Cells.Item(37, 22).Value = 4

$ws.Cells.Item(38, 5).Value = 37
$ws.Cells.Item(38, 6).Value = 280
$ws.Cells.Item(38, 7).Value = "bedrooms"
$ws.Cells.Item(38, 8).Value = "living_rooms"
$ws.Cells.Item(38, 9).Value = "distractor"
$ws.Cells.Item(38, 11).Value = "f"
$ws.Cells.Item(38, 12).Value = "stimuli/img_5mw7y.png"
$ws.Cells.Item(38, 13).Value = 72.6590909090909
$ws.Cells.Item(38, 14).Value = 50.86363636363637
$ws.Cells.Item(38, 15).Value = 61.76136363636364
$ws.Cells.Item(38, 16).Value = 44
$ws.Cells.Item(38, 17).Value = 6
$ws.Cells.Item(38, 18).Value = 6
$ws.Cells.Item(38, 19).Value = 6
$ws.Cells.Item(38, 20).Value = 6
$ws.Cells.Item(38, 21).Value = 6
$ws.Cells.Item(38, 22).Value = 6

$ws.Cells.Item(39, 5).Value = 38
$ws.Cells.Item(39, 6).Value = 281
$ws.Cells.Item(39, 7).Value = "bedrooms"
$ws.Cells.Item(39, 8).Value = "living_rooms"
$ws.Cells.Item(39, 9).Value = "distractor"
$ws.Cells.Item(39, 11).Value = "f"
$ws.Cells.Item(39, 12).Value = "stimuli/img_kost0.png"
$ws.Cells.Item(39, 13).Value = 63.09090909090909
$ws.Cells.Item(39, 14).Value = 42.77272727272727
$ws.Cells.Item(39, 15).Value = 52.93181818181819
$ws.Cells.Item(39, 16).Value = 44
$ws.Cells.Item(39, 17).Value = 5
$ws.Cells.Item(39, 18).Value = 5
$ws.Cells.Item(39, 19).Value = 5
$ws.Cells.Item(39, 20).Value = 5
$ws.Cells.Item(39, 21).Value = 5
$ws.Cells.Item(39, 22).Value = 5

$ws.Cells.Item(40, 5).Value = 39
$ws.Cells.Item(40, 6).Value = 282
$ws.Cells.Item(40, 7).Value = "bedrooms"
$ws.Cells.Item(40, 8).Value = "bedrooms"
$ws.Cells.Item(40, 9).Value = "target"
$ws.Cells.Item(40, 11).Value = "j"
$ws.Cells.Item(40, 12).Value = "stimuli/img_3bxjb.png"
$ws.Cells.Item(40, 13).Value = 87.28571428571429
$ws.Cells.Item(40, 14).Value = 72.65714285714286
$ws.Cells.Item(40, 15).Value = 79.97142857142858
$ws.Cells.Item(40, 16).Value = 35
$ws.Cells.Item(40, 17).Value = 10
$ws.Cells.Item(40, 18).Value = 10
$ws.Cells.Item(40, 19).Value = 10
$ws.Cells.Item(40, 20).Value = 9
$ws.Cells.Item(40, 21).Value = 9
$ws.Cells.Item(40, 22).Value = 10

$ws.Cells.Item(41, 5).Value = 40
$ws.Cells.Item(41, 6).Value = 283
$ws.Cells.Item(41, 7).Value = "bedrooms"
$ws.Cells.Item(41, 8).Value = "bedrooms"
$ws.Cells.Item(41, 9).Value = "target"
$ws.Cells.Item(41, 11).Value = "j"
$ws.Cells.Item(41, 12).Value = "stimuli/img_cgdyc.png"
$ws.Cells.Item(41, 13).Value = 32.93023255813954
$ws.Cells.Item(41, 14).Value = 14.04651162790698
$ws.Cells.Item(41, 15).Value = 23.48837209302326
$ws.Cells.Item(41, 16).Value = 43
$ws.Cells.Item(41, 17).Value = 1
$ws.Cells.Item(41, 18).Value = 1
$ws.Cells.Item(41, 19).Value = 1
$ws.Cells.Item(41, 20).Value = 1
$ws.Cells.Item(41, 21).Value = 1
$ws.Cells.Item(41, 22).Value = 1

# Fill in constant columns A-D for newly added rows 28-41
$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "categorization"
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 4).Value = 2

$ws.Cells.Item(29, 1).Value = 1
$ws.Cells.Item(29, 2).Value = "categorization"
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 2

$ws.Cells.Item(30, 1).Value = 1
$ws.Cells.Item(30, 2).Value = "categorization"
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 4).Value = 2

$ws.Cells.Item(31, 1).Value = 1
$ws.Cells.Item(31, 2).Value = "categorization"
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 2

$ws.Cells.Item(32, 1).Value = 1
$ws.Cells.Item(32, 2).Value = "categorization"
$ws.Cells.Item(32, 3).Value = 1
$ws.Cells.Item(32, 4).Value = 2

$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = "categorization"
$ws.Cells.Item(33, 3).Value = 1
$ws.Cells.Item(33, 4).Value = 2

$ws.Cells.Item(34, 1).Value = 1
$ws.Cells.Item(34, 2).Value = "categorization"
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 2

$ws.Cells.Item(35, 1).Value = 1
$ws.Cells.Item(35, 2).Value = "categorization"
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 4).Value = 2

$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = "categorization"
$ws.Cells.Item(36, 3).Value = 1
$ws.Cells.Item(36, 4).Value = 2

$ws.Cells.Item(37, 1).Value = 1
$ws.Cells.Item(37, 2).Value = "categorization"
$ws.Cells.Item(37, 3).Value = 1
$ws.Cells.Item(37, 4).Value = 2

$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(38, 2).Value = "categorization"
$ws.Cells.Item(38, 3).Value = 1
$ws.Cells.Item(38, 4).Value = 2

$ws.Cells.Item(39, 1).Value = 1
$ws.Cells.Item(39, 2).Value = "categorization"
$ws.Cells.Item(39, 3).Value = 1
$ws.Cells.Item(39, 4).Value = 2

$ws.Cells.Item(40, 1).Value = 1
$ws.Cells.Item(40, 2).Value = "categorization"
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(40, 4).Value = 2

$ws.Cells.Item(41, 1).Value = 1
$ws.Cells.Item(41, 2).Value = "categorization"
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(41, 4).Value = 2
